$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Interface_Category")

# Update existing cell text values
$ws.Range("C4").Value = "Requests for visual presentation of data or suggesting ways to do that"
$ws.Range("B5").Value = "Specific Request"
$ws.Range("B6").Value = "I Want Data"
$ws.Range("B7").Value = "Share Information"
$ws.Range("B8").Value = "A Nice API"
$ws.Range("C8").Value = "Requests for a programming interface for access to model output"
$ws.Range("C10").Value = "Suggestions for documentation, tutorials or other ways to help users"

# Add description for "Other" row
$ws.Range("C11").Value = "Anything Else"

# Add new row 12 "Examples"
$ws.Range("A12").Value = 12
$ws.Range("B12").Value = "Examples"
$ws.Range("C12").Value = "Examples or suggestions for places to look for good examples"

# Update the defined name range to cover the new row
$wb.Names("Interface_Category").RefersTo = "='Interface_Category'!`$A`$1:`$C`$12"
